$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D toggles (values filled in / cleared to blank) ---
$ws.Range("D6").Value = -14.2
$ws.Range("D8").ClearContents()
$ws.Range("D12").Value = -14.1
$ws.Range("D14").ClearContents()
$ws.Range("D17").Value = -14.7
$ws.Range("D18").Value = -15.2
$ws.Range("D19").ClearContents()
$ws.Range("D20").ClearContents()
$ws.Range("D23").Value = -13.9

# --- Remove the "RM 232" row (originally row 26) entirely ---
$ws.Rows.Item(26).Delete()

# --- Remove the "SC 92" row (originally row 28, now row 27 after the above delete) entirely ---
$ws.Rows.Item(27).Delete()

# --- Remaining per-cell edits on the shifted rows ---
# SC 101 is now row 27
$ws.Range("B27").Value = -20.4
$ws.Range("D27").ClearContents()

# SC 119 is now row 29
$ws.Range("B29").ClearContents()

# SC 193 is now row 32
$ws.Range("B32").ClearContents()
